$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-13, columns G-T (except unchanged cells)

# Row 2
$ws.Range("G2").Value = 0.009568666666666666
$ws.Range("H2").Value = 0.028706
$ws.Range("I2").Value = 0.0176680032866696
$ws.Range("J2").Value = 0.0176680032866696
$ws.Range("M2").Value = 29.04796866666666
$ws.Range("N2").Value = 87.14390599999999
$ws.Range("O2").Value = 0.2371972210028098
$ws.Range("P2").Value = 0.2371972210028099
$ws.Range("Q2").Value = 0.2779503295151111
$ws.Range("R2").Value = 2.501552965635999
$ws.Range("S2").Value = 0.00419080128026654
$ws.Range("T2").Value = 0.004190801280266541

# Row 3
$ws.Range("G3").Value = 0.009568666666666666
$ws.Range("H3").Value = 0.028706
$ws.Range("I3").Value = 0.0176680032866696
$ws.Range("J3").Value = 0.0176680032866696
$ws.Range("O3").Value = 0.2597953978506987
$ws.Range("P3").Value = 0.2597953978506987
$ws.Range("Q3").Value = 0.3044311233235555
$ws.Range("R3").Value = 2.739880109912
$ws.Range("S3").Value = 0.00459006594308778
$ws.Range("T3").Value = 0.004590065943087781

# Row 4
$ws.Range("G4").Value = 0.009568666666666666
$ws.Range("H4").Value = 0.028706
$ws.Range("I4").Value = 0.0176680032866696
$ws.Range("J4").Value = 0.0176680032866696
$ws.Range("M4").Value = 12.70280433333333
$ws.Range("N4").Value = 38.108413
$ws.Range("O4").Value = 0.1037273869778955
$ws.Range("P4").Value = 0.1037273869778955
$ws.Range("Q4").Value = 0.1215489003975555
$ws.Range("R4").Value = 1.093940103578
$ws.Range("S4").Value = 0.001832655814043108
$ws.Range("T4").Value = 0.001832655814043108

# Row 5
$ws.Range("G5").Value = 0.009568666666666666
$ws.Range("H5").Value = 0.028706
$ws.Range("I5").Value = 0.0176680032866696
$ws.Range("J5").Value = 0.0176680032866696
$ws.Range("M5").Value = 48.89716966666666
$ws.Range("N5").Value = 146.691509
$ws.Range("O5").Value = 0.3992799941685959
$ws.Range("P5").Value = 0.399279994168596
$ws.Range("Q5").Value = 0.4678807174837777
$ws.Range("R5").Value = 4.210926457354
$ws.Range("S5").Value = 0.007054480249272172
$ws.Range("T5").Value = 0.007054480249272174

# Row 6
$ws.Range("I6").Value = 0.9474452914149606
$ws.Range("J6").Value = 0.9474452914149606
$ws.Range("M6").Value = 29.04796866666666
$ws.Range("N6").Value = 87.14390599999999
$ws.Range("O6").Value = 0.2371972210028098
$ws.Range("P6").Value = 0.2371972210028099
$ws.Range("Q6").Value = 14.90506463427133
$ws.Range("R6").Value = 134.145581708442
$ws.Range("S6").Value = 0.224731390175826
$ws.Range("T6").Value = 0.224731390175826

# Row 7
$ws.Range("I7").Value = 0.9474452914149606
$ws.Range("J7").Value = 0.9474452914149606
$ws.Range("O7").Value = 0.2597953978506987
$ws.Range("P7").Value = 0.2597953978506987
$ws.Range("S7").Value = 0.2461419264249208
$ws.Range("T7").Value = 0.2461419264249208

# Row 8
$ws.Range("I8").Value = 0.9474452914149606
$ws.Range("J8").Value = 0.9474452914149606
$ws.Range("M8").Value = 12.70280433333333
$ws.Range("N8").Value = 38.108413
$ws.Range("O8").Value = 0.1037273869778955
$ws.Range("P8").Value = 0.1037273869778955
$ws.Range("Q8").Value = 6.518050256715666
$ws.Range("R8").Value = 58.66245231044099
$ws.Range("S8").Value = 0.09827602438298462
$ws.Range("T8").Value = 0.09827602438298463

# Row 9
$ws.Range("I9").Value = 0.9474452914149606
$ws.Range("J9").Value = 0.9474452914149606
$ws.Range("M9").Value = 48.89716966666666
$ws.Range("N9").Value = 146.691509
$ws.Range("O9").Value = 0.3992799941685959
$ws.Range("P9").Value = 0.399279994168596
$ws.Range("Q9").Value = 25.09006680219033
$ws.Range("R9").Value = 225.810601219713
$ws.Range("S9").Value = 0.3782959504312292
$ws.Range("T9").Value = 0.3782959504312292

# Row 10
$ws.Range("G10").Value = 0.018894
$ws.Range("H10").Value = 0.056682
$ws.Range("I10").Value = 0.0348867052983699
$ws.Range("J10").Value = 0.03488670529836991
$ws.Range("M10").Value = 29.04796866666666
$ws.Range("N10").Value = 87.14390599999999
$ws.Range("O10").Value = 0.2371972210028098
$ws.Range("P10").Value = 0.2371972210028099
$ws.Range("Q10").Value = 0.548832319988
$ws.Range("R10").Value = 4.939490879891999
$ws.Range("S10").Value = 0.008275029546717343
$ws.Range("T10").Value = 0.008275029546717346

# Row 11
$ws.Range("G11").Value = 0.018894
$ws.Range("H11").Value = 0.056682
$ws.Range("I11").Value = 0.0348867052983699
$ws.Range("J11").Value = 0.03488670529836991
$ws.Range("O11").Value = 0.2597953978506987
$ws.Range("P11").Value = 0.2597953978506987
$ws.Range("Q11").Value = 0.6011204950960001
$ws.Range("R11").Value = 5.410084455864
$ws.Range("S11").Value = 0.009063405482690086
$ws.Range("T11").Value = 0.009063405482690088

# Row 12
$ws.Range("G12").Value = 0.018894
$ws.Range("H12").Value = 0.056682
$ws.Range("I12").Value = 0.0348867052983699
$ws.Range("J12").Value = 0.03488670529836991
$ws.Range("M12").Value = 12.70280433333333
$ws.Range("N12").Value = 38.108413
$ws.Range("O12").Value = 0.1037273869778955
$ws.Range("P12").Value = 0.1037273869778955
$ws.Range("Q12").Value = 0.240006785074
$ws.Range("R12").Value = 2.160061065666
$ws.Range("S12").Value = 0.003618706780867813
$ws.Range("T12").Value = 0.003618706780867815

# Row 13
$ws.Range("G13").Value = 0.018894
$ws.Range("H13").Value = 0.056682
$ws.Range("I13").Value = 0.0348867052983699
$ws.Range("J13").Value = 0.03488670529836991
$ws.Range("M13").Value = 48.89716966666666
$ws.Range("N13").Value = 146.691509
$ws.Range("O13").Value = 0.3992799941685959
$ws.Range("P13").Value = 0.399279994168596
$ws.Range("Q13").Value = 0.923863123682
$ws.Range("R13").Value = 8.314768113138001
$ws.Range("S13").Value = 0.01392956348809466
$ws.Range("T13").Value = 0.01392956348809466
